# Commit: "added ability to read into Map"
# Populate the new E/F columns on rows 11-12 with a second name/value pair,
# and move the active selection to the newly-populated cell E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = "Nam1"
$ws.Range("F11").Value = "Nam2"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 2

$ws.Range("E11").Select()
